# "tela inicial pre feita"
# Replace the numeric TR-capacity values in column E (rows 5-10) of the
# "fancoletes" sheet with their formatted text-label equivalents, and
# leave the selection parked on the first changed cell (E5), matching
# the state the workbook was left in after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value  = "0,75 TR"
$ws.Range("E6").Value  = "1,00 TR"
$ws.Range("E7").Value  = "1,50 TR"
$ws.Range("E8").Value  = "2,00 TR"
$ws.Range("E9").Value  = "2,50 TR"
$ws.Range("E10").Value = "3,00 TR"

$ws.Range("E5").Select()
